# The document contains two "<id>...</id>" placeholders. Each one was
# originally split across three separate runs:
#   <id>        (Courier New, color 7f6000, sz 18)
#   p088v_aN    (default font, color 000000)
#   </id>       (Courier New, color 7f6000, sz 18)
#
# The edit merges each triple of runs into a single run (keeping the
# formatting of the first/"<id>" run) whose text is the full tag with a new
# id value: "<id>p088v_1</id>" and "<id>p088v_2</id>" respectively.
#
# Find.Execute locates a Range spanning all the runs that contain the
# matched text. Assigning Range.Text merges the matched runs into one run
# using the formatting of the first matched run - exactly reproducing the
# edit shown in the diff.

$d = $word.ActiveDocument

$rng1 = $d.Content
$null = $rng1.Find.Execute("<id>p088v_a1</id>", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
if ($rng1.Find.Found) {
    $rng1.Text = "<id>p088v_1</id>"
}

$rng2 = $d.Content
$null = $rng2.Find.Execute("<id>p088v_a2</id>", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $rng2.Text = "<id>p088v_2</id>"
}
